$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data permutation for rows 19-30 (rows 19/20 swap; 21/25/26 rotate; 22/24/27/28/29/30 rotate) ---
# Columns A,B,E,F,G,H,Q,R (and D for rows 22/30) take on the values shown below.

$ws.Range("A19").Value = 111670510.0
$ws.Range("B19").Value = 96346.0
$ws.Range("E19").Value = 620.0
$ws.Range("F19").Value = "Skogsfru"
$ws.Range("G19").Value = "Epipogium aphyllum"
$ws.Range("H19").Value = "Sw."
$ws.Range("Q19").Value = 558124.4538526792
$ws.Range("R19").Value = 7067994.321708324
$ws.Range("A20").Value = 111671179.0
$ws.Range("B20").Value = 78578.0
$ws.Range("E20").Value = 6458.0
$ws.Range("F20").Value = "Lunglav"
$ws.Range("G20").Value = "Lobaria pulmonaria"
$ws.Range("H20").Value = "(L.) Hoffm."
$ws.Range("Q20").Value = 558215.9656782644
$ws.Range("R20").Value = 7067867.520903144
$ws.Range("A21").Value = 111671226.0
$ws.Range("B21").Value = 78579.0
$ws.Range("E21").Value = 2081.0
$ws.Range("F21").Value = "Skrovellav"
$ws.Range("G21").Value = "Lobaria scrobiculata"
$ws.Range("H21").Value = "(Scop.) DC."
$ws.Range("Q21").Value = 558118.4535210516
$ws.Range("R21").Value = 7067742.103054954
$ws.Range("A22").Value = 111670477.0
$ws.Range("B22").Value = 96346.0
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 620.0
$ws.Range("F22").Value = "Skogsfru"
$ws.Range("G22").Value = "Epipogium aphyllum"
$ws.Range("H22").Value = "Sw."
$ws.Range("Q22").Value = 558155.0815836267
$ws.Range("R22").Value = 7068017.481975557
$ws.Range("A24").Value = 111670558.0
$ws.Range("Q24").Value = 558133.6011735104
$ws.Range("R24").Value = 7067979.426396712
$ws.Range("A25").Value = 111671197.0
$ws.Range("B25").Value = 78578.0
$ws.Range("E25").Value = 6458.0
$ws.Range("F25").Value = "Lunglav"
$ws.Range("G25").Value = "Lobaria pulmonaria"
$ws.Range("H25").Value = "(L.) Hoffm."
$ws.Range("Q25").Value = 558250.1783714101
$ws.Range("R25").Value = 7067936.828089682
$ws.Range("A26").Value = 111670567.0
$ws.Range("B26").Value = 96346.0
$ws.Range("E26").Value = 620.0
$ws.Range("F26").Value = "Skogsfru"
$ws.Range("G26").Value = "Epipogium aphyllum"
$ws.Range("H26").Value = "Sw."
$ws.Range("Q26").Value = 558129.9933989302
$ws.Range("R26").Value = 7067958.536170656
$ws.Range("A27").Value = 111671201.0
$ws.Range("B27").Value = 78579.0
$ws.Range("E27").Value = 2081.0
$ws.Range("F27").Value = "Skrovellav"
$ws.Range("G27").Value = "Lobaria scrobiculata"
$ws.Range("H27").Value = "(Scop.) DC."
$ws.Range("Q27").Value = 558250.1783714101
$ws.Range("R27").Value = 7067936.828089682
$ws.Range("A28").Value = 111671294.0
$ws.Range("B28").Value = 78578.0
$ws.Range("E28").Value = 6458.0
$ws.Range("F28").Value = "Lunglav"
$ws.Range("G28").Value = "Lobaria pulmonaria"
$ws.Range("H28").Value = "(L.) Hoffm."
$ws.Range("Q28").Value = 558118.4535210516
$ws.Range("R28").Value = 7067742.103054954
$ws.Range("A29").Value = 111670497.0
$ws.Range("B29").Value = 96346.0
$ws.Range("E29").Value = 620.0
$ws.Range("F29").Value = "Skogsfru"
$ws.Range("G29").Value = "Epipogium aphyllum"
$ws.Range("H29").Value = "Sw."
$ws.Range("Q29").Value = 558159.8619213518
$ws.Range("R29").Value = 7068022.886732788
$ws.Range("A30").Value = 111671188.0
$ws.Range("B30").Value = 78605.0
$ws.Range("D30").Value = "LC"
$ws.Range("E30").Value = 6462.0
$ws.Range("F30").Value = "Stuplav"
$ws.Range("G30").Value = "Nephroma bellum"
$ws.Range("H30").Value = "(Spreng.) Tuck."
$ws.Range("Q30").Value = 558215.9329796816
$ws.Range("R30").Value = 7067869.292590594

# --- "Kön" (L) column: present (empty) only for Epipogium aphyllum / Skogsfru rows ---
# Rows gaining an (empty) L cell:
$ws.Range("L19").Font.Bold = $false
$ws.Range("L22").Font.Bold = $false
$ws.Range("L26").Font.Bold = $false
$ws.Range("L29").Font.Bold = $false
# Rows losing their (empty) L cell:
$ws.Range("L20").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("L27").ClearContents()
$ws.Range("L30").ClearContents()
